# Appends 13 new flight-arrival rows (Friday, Jan 13) to the "Main Data" sheet,
# extending the table from row 229 to row 242 (NUMBER 229..241).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 230
$ws.Cells.Item(230, 1).Value = 229
$ws.Cells.Item(230, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(230, 3).Value = "2:40 PM"
$ws.Cells.Item(230, 4).Value = "LO3883"
$ws.Cells.Item(230, 5).Value = "Warsaw"
$ws.Cells.Item(230, 6).Value = "(WAW)"
$ws.Cells.Item(230, 7).Value = "LOT "
$ws.Cells.Item(230, 8).Value = "E170"
$ws.Cells.Item(230, 9).Value = "(SP-LDH)"
$ws.Cells.Item(230, 10).Value = "2:36 PM"
$ws.Cells.Item(230, 11).Font.Size = 11
$ws.Cells.Item(230, 12).Value = "0 hours, -4 minutes"
$ws.Cells.Item(230, 13).Font.Size = 11

# Row 231
$ws.Cells.Item(231, 1).Value = 230
$ws.Cells.Item(231, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(231, 3).Value = "3:20 PM"
$ws.Cells.Item(231, 4).Value = "FR2465"
$ws.Cells.Item(231, 5).Value = "London"
$ws.Cells.Item(231, 6).Value = "(STN)"
$ws.Cells.Item(231, 7).Value = "Ryanair "
$ws.Cells.Item(231, 8).Value = "B738"
$ws.Cells.Item(231, 9).Value = "(SP-RKB)"
$ws.Cells.Item(231, 10).Value = "3:42 PM"
$ws.Cells.Item(231, 11).Font.Size = 11
$ws.Cells.Item(231, 12).Value = "0 hours, 22 minutes"
$ws.Cells.Item(231, 13).Font.Size = 11

# Row 232
$ws.Cells.Item(232, 1).Value = 231
$ws.Cells.Item(232, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(232, 3).Value = "4:10 PM"
$ws.Cells.Item(232, 4).Value = "KL1815"
$ws.Cells.Item(232, 5).Value = "Amsterdam"
$ws.Cells.Item(232, 6).Value = "(AMS)"
$ws.Cells.Item(232, 7).Value = "KLM "
$ws.Cells.Item(232, 8).Value = "E295"
$ws.Cells.Item(232, 9).Value = "(PH-NXB)"
$ws.Cells.Item(232, 10).Value = "3:54 PM"
$ws.Cells.Item(232, 11).Font.Size = 11
$ws.Cells.Item(232, 12).Value = "0 hours, -16 minutes"
$ws.Cells.Item(232, 13).Font.Size = 11

# Row 233
$ws.Cells.Item(233, 1).Value = 232
$ws.Cells.Item(233, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(233, 3).Value = "4:20 PM"
$ws.Cells.Item(233, 4).Value = "FR6893"
$ws.Cells.Item(233, 5).Value = "Dortmund"
$ws.Cells.Item(233, 6).Value = "(DTM)"
$ws.Cells.Item(233, 7).Value = "Ryanair "
$ws.Cells.Item(233, 8).Value = "B738"
$ws.Cells.Item(233, 9).Value = "(SP-RSB)"
$ws.Cells.Item(233, 10).Value = "4:03 PM"
$ws.Cells.Item(233, 11).Font.Size = 11
$ws.Cells.Item(233, 12).Value = "0 hours, -17 minutes"
$ws.Cells.Item(233, 13).Font.Size = 11

# Row 234
$ws.Cells.Item(234, 1).Value = 233
$ws.Cells.Item(234, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(234, 3).Value = "5:10 PM"
$ws.Cells.Item(234, 4).Value = "W61072"
$ws.Cells.Item(234, 5).Value = "Eindhoven"
$ws.Cells.Item(234, 6).Value = "(EIN)"
$ws.Cells.Item(234, 7).Value = "Wizz Air "
$ws.Cells.Item(234, 8).Value = "A321"
$ws.Cells.Item(234, 9).Value = "(HA-LXD)"
$ws.Cells.Item(234, 10).Value = "5:04 PM"
$ws.Cells.Item(234, 11).Font.Size = 11
$ws.Cells.Item(234, 12).Value = "0 hours, -6 minutes"
$ws.Cells.Item(234, 13).Font.Size = 11

# Row 235
$ws.Cells.Item(235, 1).Value = 234
$ws.Cells.Item(235, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(235, 3).Value = "5:40 PM"
$ws.Cells.Item(235, 4).Value = "FR5669"
$ws.Cells.Item(235, 5).Value = "Edinburgh"
$ws.Cells.Item(235, 6).Value = "(EDI)"
$ws.Cells.Item(235, 7).Value = "Ryanair "
$ws.Cells.Item(235, 8).Value = "B38M"
$ws.Cells.Item(235, 9).Value = "(EI-HES)"
$ws.Cells.Item(235, 10).Value = "5:30 PM"
$ws.Cells.Item(235, 11).Font.Size = 11
$ws.Cells.Item(235, 12).Value = "0 hours, -10 minutes"
$ws.Cells.Item(235, 13).Font.Size = 11

# Row 236
$ws.Cells.Item(236, 1).Value = 235
$ws.Cells.Item(236, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(236, 3).Value = "6:25 PM"
$ws.Cells.Item(236, 4).Value = "W61094"
$ws.Cells.Item(236, 5).Value = "Dortmund"
$ws.Cells.Item(236, 6).Value = "(DTM)"
$ws.Cells.Item(236, 7).Value = "Wizz Air "
$ws.Cells.Item(236, 8).Value = "A320"
$ws.Cells.Item(236, 9).Value = "(HA-LWR)"
$ws.Cells.Item(236, 10).Value = "6:09 PM"
$ws.Cells.Item(236, 11).Font.Size = 11
$ws.Cells.Item(236, 12).Value = "0 hours, -16 minutes"
$ws.Cells.Item(236, 13).Font.Size = 11

# Row 237
$ws.Cells.Item(237, 1).Value = 236
$ws.Cells.Item(237, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(237, 3).Value = "6:55 PM"
$ws.Cells.Item(237, 4).Value = "LH1362"
$ws.Cells.Item(237, 5).Value = "Frankfurt"
$ws.Cells.Item(237, 6).Value = "(FRA)"
$ws.Cells.Item(237, 7).Value = "Lufthansa "
$ws.Cells.Item(237, 8).Value = "CRJ9"
$ws.Cells.Item(237, 9).Value = "(D-ACNO)"
$ws.Cells.Item(237, 10).Value = "6:46 PM"
$ws.Cells.Item(237, 11).Font.Size = 11
$ws.Cells.Item(237, 12).Value = "0 hours, -9 minutes"
$ws.Cells.Item(237, 13).Font.Size = 11

# Row 238
$ws.Cells.Item(238, 1).Value = 237
$ws.Cells.Item(238, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(238, 3).Value = "9:15 PM"
$ws.Cells.Item(238, 4).Value = "FR6392"
$ws.Cells.Item(238, 5).Value = "London"
$ws.Cells.Item(238, 6).Value = "(STN)"
$ws.Cells.Item(238, 7).Value = "Lauda Europe "
$ws.Cells.Item(238, 8).Value = "A320"
$ws.Cells.Item(238, 9).Value = "(9H-LMH)"
$ws.Cells.Item(238, 10).Value = "9:20 PM"
$ws.Cells.Item(238, 11).Font.Size = 11
$ws.Cells.Item(238, 12).Value = "0 hours, 5 minutes"
$ws.Cells.Item(238, 13).Font.Size = 11

# Row 239
$ws.Cells.Item(239, 1).Value = 238
$ws.Cells.Item(239, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(239, 3).Value = "9:20 PM"
$ws.Cells.Item(239, 4).Value = "RR7990"
$ws.Cells.Item(239, 5).Value = "Gran Canaria"
$ws.Cells.Item(239, 6).Value = "(LPA)"
$ws.Cells.Item(239, 7).Value = "Ryanair "
$ws.Cells.Item(239, 8).Value = "B738"
$ws.Cells.Item(239, 9).Value = "(SP-RSN)"
$ws.Cells.Item(239, 10).Value = "9:09 PM"
$ws.Cells.Item(239, 11).Font.Size = 11
$ws.Cells.Item(239, 12).Value = "0 hours, -11 minutes"
$ws.Cells.Item(239, 13).Font.Size = 11

# Row 240
$ws.Cells.Item(240, 1).Value = 239
$ws.Cells.Item(240, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(240, 3).Value = "9:27 PM"
$ws.Cells.Item(240, 4).Value = "3V4101"
$ws.Cells.Item(240, 5).Value = "Kaunas"
$ws.Cells.Item(240, 6).Value = "(KUN)"
$ws.Cells.Item(240, 7).Value = "ASL Airlines "
$ws.Cells.Item(240, 8).Value = "B734"
$ws.Cells.Item(240, 9).Value = "(OE-IAB)"
$ws.Cells.Item(240, 10).Value = "9:06 PM"
$ws.Cells.Item(240, 11).Font.Size = 11
$ws.Cells.Item(240, 12).Value = "0 hours, -21 minutes"
$ws.Cells.Item(240, 13).Font.Size = 11

# Row 241
$ws.Cells.Item(241, 1).Value = 240
$ws.Cells.Item(241, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(241, 3).Value = "10:00 PM"
$ws.Cells.Item(241, 4).Value = "FR6388"
$ws.Cells.Item(241, 5).Value = "Athens"
$ws.Cells.Item(241, 6).Value = "(ATH)"
$ws.Cells.Item(241, 7).Value = "Ryanair "
$ws.Cells.Item(241, 8).Value = "B738"
$ws.Cells.Item(241, 9).Value = "(SP-RSB)"
$ws.Cells.Item(241, 10).Value = "9:47 PM"
$ws.Cells.Item(241, 11).Font.Size = 11
$ws.Cells.Item(241, 12).Value = "0 hours, -13 minutes"
$ws.Cells.Item(241, 13).Font.Size = 11

# Row 242
$ws.Cells.Item(242, 1).Value = 241
$ws.Cells.Item(242, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(242, 3).Value = "10:20 PM"
$ws.Cells.Item(242, 4).Value = "3Z7607"
$ws.Cells.Item(242, 5).Value = "Salalah"
$ws.Cells.Item(242, 6).Value = "(SLL)"
$ws.Cells.Item(242, 7).Value = "Smartwings "
$ws.Cells.Item(242, 8).Value = "B38M"
$ws.Cells.Item(242, 9).Value = "(OK-SWC)"
$ws.Cells.Item(242, 10).Value = "10:21 PM"
$ws.Cells.Item(242, 11).Font.Size = 11
$ws.Cells.Item(242, 12).Value = "0 hours, 1 minutes"
$ws.Cells.Item(242, 13).Font.Size = 11

